$wb = $excel.ActiveWorkbook

# --- Sheet "pomn": originally columns A..H, where H only held a header
#     ("f") with no data beneath it (A1:H6 dimension, H2:H6 empty).
#     Shift A..G right into B..H (overwriting the stray, data-less H),
#     then fill in a new leading id column A and relabel B1.
$ws1 = $wb.Worksheets.Item("pomn")
for ($r = 1; $r -le 6; $r++) {
    $ws1.Cells.Item($r, 8).Value = $ws1.Cells.Item($r, 7).Value2
    $ws1.Cells.Item($r, 7).Value = $ws1.Cells.Item($r, 6).Value2
    $ws1.Cells.Item($r, 6).Value = $ws1.Cells.Item($r, 5).Value2
    $ws1.Cells.Item($r, 5).Value = $ws1.Cells.Item($r, 4).Value2
    $ws1.Cells.Item($r, 4).Value = $ws1.Cells.Item($r, 3).Value2
    $ws1.Cells.Item($r, 3).Value = $ws1.Cells.Item($r, 2).Value2
    $ws1.Cells.Item($r, 2).Value = $ws1.Cells.Item($r, 1).Value2
}
$ws1.Range("B1").Value = "l"
$ws1.Range("A1").Value = "lp"
$ws1.Range("A2").Value = 1
$ws1.Range("A3").Value = 2
$ws1.Range("A4").Value = 3
$ws1.Range("A5").Value = 4
$ws1.Range("A6").Value = 5
$ws1.Range("A7").Select()

# --- Sheet "pow": originally columns A..G (no stray trailing column).
#     Shift A..G right into B..H, then fill in the new leading id
#     column A and relabel B1.
$ws2 = $wb.Worksheets.Item("pow")
for ($r = 1; $r -le 6; $r++) {
    $ws2.Cells.Item($r, 8).Value = $ws2.Cells.Item($r, 7).Value2
    $ws2.Cells.Item($r, 7).Value = $ws2.Cells.Item($r, 6).Value2
    $ws2.Cells.Item($r, 6).Value = $ws2.Cells.Item($r, 5).Value2
    $ws2.Cells.Item($r, 5).Value = $ws2.Cells.Item($r, 4).Value2
    $ws2.Cells.Item($r, 4).Value = $ws2.Cells.Item($r, 3).Value2
    $ws2.Cells.Item($r, 3).Value = $ws2.Cells.Item($r, 2).Value2
    $ws2.Cells.Item($r, 2).Value = $ws2.Cells.Item($r, 1).Value2
}
$ws2.Range("B1").Value = "l"
$ws2.Range("A1").Value = "lp"
$ws2.Range("A2").Value = 1
$ws2.Range("A3").Value = 2
$ws2.Range("A4").Value = 3
$ws2.Range("A5").Value = 4
$ws2.Range("A6").Value = 5
$ws2.Range("A7").Select()

$ws2.Activate()
